# LMS-2340 Finished implementing changes to templates and handlers.
#
# The "openbis-metadata" sheet's Header-Format row (row 2) gets its Value
# cell (B2) filled in with a sample "/TEST/TEST/TEST" path, and the
# selection cursor ends up one cell to the right of where it previously
# was (B9 -> C9).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openbis-metadata")

# Fill in B2 (previously blank) with the sample value; existing cell
# style (s="5") is left untouched by a plain .Value assignment.
$ws.Range("B2").Value = "/TEST/TEST/TEST"

# Recorded cursor/selection move on this sheet.
$ws.Activate()
$ws.Range("C9").Select()
